$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trading History")

# Insert a new row at position 5, pushing existing rows 5-8 down to 6-9
$ws.Rows("5:5").Insert()

# The inserted row inherits the header row's formatting/shape; strip it
# back down to the plain (unstyled) look used by the other data rows,
# and drop the cells that don't actually hold data in this sheet.
$ws.Range("A5:N5").ClearFormats()
$ws.Range("W5:AB5").ClearFormats()
$ws.Range("H5").Value = $null
$ws.Range("K5:N5").Value = $null
$ws.Range("W5:AB5").Value = $null
$ws.Range("A5").NumberFormat = "yyyy-mm-dd h:mm:ss"

# Populate the newly inserted row 5 with the new trade entry
$ws.Range("A5").Value = 46062
$ws.Range("B5").Value = "NSE"
$ws.Range("C5").Value = "Buy"
$ws.Range("D5").Value = 100
$ws.Range("E5").Value = 36.7
$ws.Range("F5").Value = 3688
$ws.Range("G5").Value = "CN#252611665409"
$ws.Range("I5").Value = 18
$ws.Range("J5").Formula = "=Index!`$C`$2"
